$wb = $excel.ActiveWorkbook

# OFF sheet - Row 2 (H)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 286
$wsOff.Range("C2").Value = 221
$wsOff.Range("D2").Value = 72
$wsOff.Range("E2").Value = 32
$wsOff.Range("G2").Value = 4

# DEF sheet - Row 2 (H)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 420
$wsDef.Range("C2").Value = 293
$wsDef.Range("D2").Value = 102
$wsDef.Range("E2").Value = 44
$wsDef.Range("G2").Value = 7
